$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header text: unit changes from GWh to MWh
$ws.Range("B1").Value = "Electricity Generation (MWh) from Solar Photovoltaics"

# Move the active selection to E3 (matches the saved selection state in the file)
$ws.Range("E3").Select()
